# Q30644134-en.xlsx — "added one json for time bucket analysis"
#
# Row 2 and row 3 swap their "title" (column A) text and their "uri"
# (column E) text/hyperlink:
#   - A2: "U. S. Electoral College"                       -> "United States presidential election of 1960"
#   - A3: "United States presidential election of 1960"    -> "U. S. Electoral College"
#   - E2: archives.gov link (...html#1960)                 -> britannica.com link (no fragment)
#   - E3: britannica.com link (no fragment)                 -> archives.gov link (...html#1960)
#
# The hyperlink "location" (the #1960 in-doc fragment) moves from E2 to E3
# along with the new target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the displayed titles in column A -------------------------------
$ws.Range("A2").Value = "United States presidential election of 1960"
$ws.Range("A3").Value = "U. S. Electoral College"

# remember the current "Hyperlink" cell style so the rebuilt hyperlinks in
# column E keep looking the same (blue/underlined) instead of picking up a
# freshly minted (but equivalent) style index
$linkStyle = $ws.Range("E2").Style

# --- rebuild the hyperlinks in column E so they point at the swapped URLs
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.britannica.com/event/United-States-presidential-election-of-1960")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.archives.gov/federal-register/electoral-college/votes/1953_1957.html", "1960")

# the display text for each uri cell should match its new target
$ws.Range("E2").Value = "https://www.britannica.com/event/United-States-presidential-election-of-1960"
$ws.Range("E3").Value = "https://www.archives.gov/federal-register/electoral-college/votes/1953_1957.html#1960"

# restore the original hyperlink styling
$ws.Range("E2").Style = $linkStyle
$ws.Range("E3").Style = $linkStyle
